# Fruta / hortaliza, semanal
# Update weekly price data for Pomelo (rows 2-7 and 10-14).
# Columns touched: D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado),
# Q (Unidad de comercializacion), S (Precio $/Kg)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ D = 44216; M = 55; N = 11000; O = 12000; P = 11545; Q = "`$/caja 14 kilos empedrada"; S = 825 }
    3  = @{ D = 44253; M = 90; N = 12000; O = 13000; P = 12667; Q = "`$/caja 14 kilos empedrada"; S = 905 }
    4  = @{ D = 45152; M = 60; N = 16000; O = 16000; P = 16000; Q = "`$/caja 14 kilos empedrada"; S = 1143 }
    5  = @{ D = 44181; M = 65; N = 9000;  O = 10000; P = 9462;  Q = "`$/caja 14 kilos empedrada"; S = 676 }
    6  = @{ D = 44172; M = 90; N = 8500;  O = 9000;  P = 8806;  Q = "`$/caja 14 kilos empedrada"; S = 629 }
    7  = @{ D = 44210; M = 70; N = 10000; O = 11000; P = 10357; Q = "`$/caja 14 kilos empedrada"; S = 740 }
    10 = @{ D = 45155; M = 60; N = 15000; O = 15000; P = 15000; Q = "`$/caja 14 kilos empedrada"; S = 1071 }
    11 = @{ D = 45142; M = 30; N = 15000; O = 15000; P = 15000; Q = "`$/caja 14 kilos empedrada"; S = 1071 }
    12 = @{ D = 45142; M = 30; N = 14000; O = 14000; P = 14000; Q = "`$/caja 14 kilos granel"; S = 1000 }
    13 = @{ D = 45138; M = 50; N = 14000; O = 14000; P = 14000; Q = "`$/caja 14 kilos granel"; S = 1000 }
    14 = @{ D = 45140; M = 30; N = 15000; O = 15000; P = 15000; Q = "`$/caja 14 kilos granel"; S = 1071 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("S$row").Value = $vals.S
}
